$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 679.1177
$ws.Range("I15").Value = 679.1177
$ws.Range("K15").Value = 2037.3531
$ws.Range("M15").Value = -1868.3531

$ws.Range("H40").Value = 8311.546
$ws.Range("I40").Value = 3425
$ws.Range("K40").Value = 3425
$ws.Range("M40").Value = -3250

$ws.Range("H98").Value = 2116.75
$ws.Range("I98").Value = 704.8570999999999
$ws.Range("K98").Value = 704.8570999999999
$ws.Range("M98").Value = 793.1429000000001

$ws.Range("H100").Value = 2520.3044
$ws.Range("I100").Value = 2522.2856
$ws.Range("K100").Value = 2522.2856
$ws.Range("M100").Value = -1981.2856

$ws.Range("H122").Value = 2116.75
$ws.Range("I122").Value = 704.8570999999999
$ws.Range("K122").Value = 2114.5713
$ws.Range("M122").Value = 335.4287000000004

$ws.Range("H137").Value = 2936.0588
$ws.Range("J137").Value = 4534.375
$ws.Range("L137").Value = 13603.125
$ws.Range("N137").Value = -18703.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8160.222
$ws.Range("I32").Value = 7169.7354
$ws.Range("K32").Value = 7169.7354
$ws.Range("M32").Value = -6882.7354

$ws.Range("H74").Value = 2938.0715
$ws.Range("I74").Value = 2938.0715
$ws.Range("K74").Value = 2938.0715
$ws.Range("M74").Value = -2064.0715

$ws.Range("H77").Value = 2938.0715
$ws.Range("I77").Value = 2938.0715
$ws.Range("K77").Value = 14690.3575
$ws.Range("M77").Value = -10322.3575

$ws.Range("H97").Value = 466.2
$ws.Range("I97").Value = 466.2
$ws.Range("K97").Value = 466.2
$ws.Range("M97").Value = 29.80000000000001

$ws.Range("H102").Value = 1252.2222
$ws.Range("I102").Value = 1252.2222
$ws.Range("K102").Value = 1252.2222
$ws.Range("M102").Value = 369.7778000000001

$ws.Range("H122").Value = 2299.25
$ws.Range("I122").Value = 2299.25
$ws.Range("K122").Value = 6897.75
$ws.Range("M122").Value = -4447.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4781.385
$ws.Range("I20").Value = 2571.111
$ws.Range("J20").Value = 9754.5
$ws.Range("K20").Value = 2571.111
$ws.Range("L20").Value = 9754.5
$ws.Range("M20").Value = -2324.111
$ws.Range("N20").Value = -10248.5

$ws.Range("H37").Value = 1339.5555
$ws.Range("I37").Value = 882
$ws.Range("K37").Value = 882
$ws.Range("M37").Value = -745

$ws.Range("H86").Value = 525
$ws.Range("I86").Value = 450
$ws.Range("J86").Value = 600
$ws.Range("K86").Value = 450
$ws.Range("L86").Value = 600
$ws.Range("M86").Value = 673
$ws.Range("N86").Value = -2846

$ws.Range("H89").Value = 525
$ws.Range("I89").Value = 450
$ws.Range("J89").Value = 600
$ws.Range("K89").Value = 2250
$ws.Range("L89").Value = 3000
$ws.Range("M89").Value = 3366
$ws.Range("N89").Value = -14232

$ws.Range("H94").Value = 1776.2
$ws.Range("J94").Value = 2641.6
$ws.Range("L94").Value = 2641.6
$ws.Range("N94").Value = -3543.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 8613.200000000001
$ws.Range("I41").Value = 3266.5
$ws.Range("K41").Value = 3266.5
$ws.Range("M41").Value = -2838.5

$ws.Range("H88").Value = 31171.5
$ws.Range("J88").Value = 31171.5
$ws.Range("L88").Value = 31171.5
$ws.Range("N88").Value = -31983.5

$ws.Range("H91").Value = 31171.5
$ws.Range("J91").Value = 31171.5
$ws.Range("L91").Value = 31171.5
$ws.Range("N91").Value = -33979.5

$ws.Range("H92").Value = 32183.666
$ws.Range("J92").Value = 32183.666
$ws.Range("L92").Value = 32183.666
$ws.Range("N92").Value = -37175.666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H57").Value = 2000
$ws.Range("I57").Value = 2000
$ws.Range("K57").Value = 6000
$ws.Range("M57").Value = -5441

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3695.4285
$ws.Range("J80").Value = 3132.6667
$ws.Range("L80").Value = 3132.6667
$ws.Range("N80").Value = -5128.6667

$ws.Range("H83").Value = 3695.4285
$ws.Range("J83").Value = 3132.6667
$ws.Range("L83").Value = 15663.3335
$ws.Range("N83").Value = -25647.3335

$ws.Range("H97").Value = 1754.5
$ws.Range("I97").Value = 1010
$ws.Range("J97").Value = 2499
$ws.Range("K97").Value = 1010
$ws.Range("L97").Value = 2499
$ws.Range("M97").Value = -514
$ws.Range("N97").Value = -3491

$ws.Range("H122").Value = 4533.421
$ws.Range("I122").Value = 4772.9414
$ws.Range("J122").Value = 2497.5
$ws.Range("K122").Value = 14318.8242
$ws.Range("L122").Value = 7492.5
$ws.Range("M122").Value = -11868.8242
$ws.Range("N122").Value = -12392.5

$ws.Range("H132").Value = 3853.889
$ws.Range("I132").Value = 2383.7144
$ws.Range("K132").Value = 7151.1432
$ws.Range("M132").Value = -4621.1432

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 774
$ws.Range("I7").Value = 774
$ws.Range("K7").Value = 774
$ws.Range("M7").Value = -662

$ws.Range("H22").Value = 2294.6
$ws.Range("I22").Value = 744.8
$ws.Range("J22").Value = 3069.5
$ws.Range("K22").Value = 744.8
$ws.Range("L22").Value = 3069.5
$ws.Range("M22").Value = -449.8
$ws.Range("N22").Value = -3659.5

$ws.Range("H27").Value = 2294.6
$ws.Range("I27").Value = 744.8
$ws.Range("J27").Value = 3069.5
$ws.Range("K27").Value = 744.8
$ws.Range("L27").Value = 3069.5
$ws.Range("M27").Value = -637.8
$ws.Range("N27").Value = -3283.5

$ws.Range("H40").Value = 17554.111
$ws.Range("I40").Value = 18997.834
$ws.Range("K40").Value = 18997.834
$ws.Range("M40").Value = -18861.834

$ws.Range("H55").Value = 1051.8572
$ws.Range("I55").Value = 307.25
$ws.Range("J55").Value = 1349.7
$ws.Range("K55").Value = 307.25
$ws.Range("L55").Value = 1349.7
$ws.Range("M55").Value = -134.25
$ws.Range("N55").Value = -1695.7

$ws.Range("H82").Value = 6898.1665
$ws.Range("I82").Value = 3890
$ws.Range("J82").Value = 7499.8
$ws.Range("K82").Value = 3890
$ws.Range("L82").Value = 7499.8
$ws.Range("M82").Value = -3529
$ws.Range("N82").Value = -8221.799999999999

$ws.Range("H85").Value = 6898.1665
$ws.Range("I85").Value = 3890
$ws.Range("J85").Value = 7499.8
$ws.Range("K85").Value = 3890
$ws.Range("L85").Value = 7499.8
$ws.Range("M85").Value = -2642
$ws.Range("N85").Value = -9995.799999999999

$ws.Range("H126").Value = 774
$ws.Range("I126").Value = 774
$ws.Range("K126").Value = 2322
$ws.Range("M126").Value = 148

$ws.Range("H132").Value = 10374.75
$ws.Range("I132").Value = 9750
$ws.Range("K132").Value = 29250
$ws.Range("M132").Value = -26720

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1182.9286
$ws.Range("I100").Value = 418
$ws.Range("K100").Value = 836
$ws.Range("M100").Value = -295

$ws.Range("H126").Value = 2061.2727
$ws.Range("I126").Value = 1630.5555
$ws.Range("K126").Value = 4891.666499999999
$ws.Range("M126").Value = -2421.666499999999

$ws.Range("H132").Value = 3384.875
$ws.Range("J132").Value = 2774.5
$ws.Range("L132").Value = 8323.5
$ws.Range("N132").Value = -13383.5
